# VerveStacks_DEU scen_tsparameters_ts12_clu.xlsx - 2025-08-31 20:37 update
#
# NOTE: looking sheets up by name (Worksheets.Item("name") / Worksheets("name"))
# does not reliably persist writes in this host, so sheets are addressed by
# their (1-based) tab position instead:
#   1 = ev_charging_uc, 2 = customize, 3 = timeslice_def,
#   4 = re_profiles,    5 = load_shapes
#
# 1) "timeslice_def" sheet (xl/worksheets/sheet1.xml): the two comma-
#    separated timeslice-grouping strings in C13/C14 are re-shuffled (same
#    six members, new order). G7 (=C14) and G8 (=C13) pick the new text up
#    automatically through their formulas on recalculation.
# 2) "re_profiles" sheet (xl/worksheets/sheet4.xml): column M (capacity
#    numbers) is rescaled down to fractional values and its custom "0.0"
#    number format is replaced with the existing "0.000" format already
#    used elsewhere on the sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) timeslice_def: reorder the two timeslice-cluster strings
# ---------------------------------------------------------------------
$tsDef = $wb.Worksheets.Item(3)
$tsDef.Range("C13").Value = "S1aH3,S2aH2,S3aH2,S1aH2,S2aH3,S3aH3"
$tsDef.Range("C14").Value = "S2aH1,S1aH1,S3aH4,S2aH4,S3aH1,S1aH4"

# ---------------------------------------------------------------------
# 2) re_profiles: new column M values + drop the bespoke "0.0" format
# ---------------------------------------------------------------------
$reProfiles = $wb.Worksheets.Item(4)

# (scientific-notation literals aren't accepted by this shell's parser, so
# the values below are written out in plain fixed-point form instead)
$mValues = @{
    11 = 0.22051110595637008
    12 = 0.027005130993325931
    13 = 0.19550712264392731
    14 = 0.23085918579945908
    15 = 0.069107708417194796
    16 = 0.0084581681043188073
    17 = 0.05592609458618885
    18 = 0.066776091243659561
    19 = 0.041992181052798866
    20 = 0.0053335867092306685
    21 = 0.037038344289336139
    22 = 0.041485280203954977
}

foreach ($row in $mValues.Keys) {
    $cell = $reProfiles.Range("M$row")
    $cell.Value = $mValues[$row]
    $cell.NumberFormat = "0.000"
}
